{"js": "// The document contains three occurrences of an `<id>...</id>` marker that\n// is currently split across three separate runs, e.g.:\n//   run1: \"<id>\"      (Courier New, color 7f6000, sz 18)\n//   run2: \"p067v_1\"   (color 000000)\n//   run3: \"</id>\"     (Courier New, color 7f6000, sz 18)\n// The edit collapses each triple into a single run containing the full\n// \"<id>p067v_N</id>\" text, keeping the formatting of the first run.\nconst ids = [\"p067v_1\", \"p067v_2\", \"p067v_3\"];\n\nfor (const id of ids) {\n  const target = `<id>${id}</id>`;\n  const results = context.document.body.search(target, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    // Replacing the whole matched range with its own text merges the\n    // underlying runs into a single run that uses the formatting of the\n    // first run in the matched range (the \"<id>\" run), exactly mirroring\n    // the OOXML collapse described by the diff.\n    range.insertText(target, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document holds three \"<id>...</id>\" markers, each currently split\n# across three separate runs:\n#   \"<id>\"      (Courier New, color 7f6000, sz 18)\n#   \"p067v_N\"   (color 000000)\n#   \"</id>\"     (Courier New, color 7f6000, sz 18)\n# Find-and-replace each full marker with itself; Word collapses the\n# matched runs into a single run (using the formatting of the first run\n# in the match), which merges the three runs into one run containing\n# \"<id>p067v_N</id>\", matching the OOXML change in the diff.\n$d = $word.ActiveDocument\n$ids = @(\"p067v_1\", \"p067v_2\", \"p067v_3\")\n\nforeach ($id in $ids) {\n    $target = \"<id>$id</id>\"\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $target\n    $find.Replacement.Text = $target\n    $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 1)\n}\n"}
